# Rename the three worksheets (data/order preserved, names updated)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "运输部"
$ws2.Name = "财务部"
$ws3.Name = "会计部"

# ---- Headers for 财务部 (sheet2) ----
$ws2.Range("A1").Value = "工号"
$ws2.Range("B1").Value = "姓名"
$ws2.Range("C1").Value = "部门"
$ws2.Range("D1").Value = "日期"

# ---- Headers for 会计部 (sheet3) ----
$ws3.Range("A1").Value = "工号"
$ws3.Range("B1").Value = "姓名"
$ws3.Range("C1").Value = "部门"
$ws3.Range("D1").Value = "日期"

# ---- Department column first (controls shared-string order) ----
$ws2.Range("C2").Value = "财务部"
$ws3.Range("C2").Value = "会计部"

# ---- Names for 会计部 (sheet3) ----
$ws3.Range("B2").Value = "李1"
$ws3.Range("B3").Value = "李2"
$ws3.Range("B4").Value = "李3"

# ---- Names for 财务部 (sheet2) ----
$ws2.Range("B2").Value = "童3"
$ws2.Range("B3").Value = "童4"
$ws2.Range("B4").Value = "童5"
$ws2.Range("B5").Value = "童6"
$ws2.Range("B6").Value = "童7"
$ws2.Range("B7").Value = "童8"
$ws2.Range("B8").Value = "童9"

# ---- Remaining columns: id numbers, department (rows 3-8), dates ----
$ws3.Range("A2").Value = 1
$ws3.Range("A3").Value = 2
$ws3.Range("A4").Value = 3

$ws3.Range("C3").Value = "会计部"
$ws3.Range("C4").Value = "会计部"

$ws3.Range("D2").Value = "2019-12-10"
$ws3.Range("D3").Value = "2019-12-11"
$ws3.Range("D4").Value = "2019-12-12"

$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Value = 6
$ws2.Range("A8").Value = 7

$ws2.Range("C3").Value = "财务部"
$ws2.Range("C4").Value = "财务部"
$ws2.Range("C5").Value = "财务部"
$ws2.Range("C6").Value = "财务部"
$ws2.Range("C7").Value = "财务部"
$ws2.Range("C8").Value = "财务部"

$ws2.Range("D2").Value = "2019-12-10"
$ws2.Range("D3").Value = "2019-12-11"
$ws2.Range("D4").Value = "2019-12-12"
$ws2.Range("D5").Value = "2019-12-13"
$ws2.Range("D6").Value = "2019-12-14"
$ws2.Range("D7").Value = "2019-12-15"
$ws2.Range("D8").Value = "2019-12-16"

# ---- Column D widths (bestFit-style autosize applied originally) ----
$ws2.Columns.Item(4).ColumnWidth = 11.125
$ws3.Columns.Item(4).ColumnWidth = 11.125

# ---- Update the selections on each sheet to match the new working cells ----
$ws2.Range("H14").Select()
$ws3.Range("H24").Select()
$ws1.Range("D27").Select()
